# FrSky Passthrough protocol: add GPS advanced fix / RTK info, replacing the
# old "Vertical dilution of precision" bit-field entry (rows 23-24).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Passthrough")

# Row 23/24 used to be one merged bit-field ("Vertical dilution of precision")
# spanning B23:B24. Break that merge since the new content is two independent
# single-row fields.
$ws.Range("B23:B24").UnMerge()

# Row 23: GPS advanced fix (2 bits) + description of the RTK states.
$ws.Range("B23").Value = "GPS advanced fix"
$ws.Range("C23").Value = "N/A"
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = "0: no advanced fix, 1: DGPS, 2: RTK_FLOAT, 3: RTK_FIXED"

# Row 24: RESERVED bits filling out the rest of the byte (6 bits).
$ws.Range("B24").Value = "RESERVED"
$ws.Range("C24").ClearContents()
$ws.Range("D24").Value = 6

# Update view state to reflect where the author ended up editing.
$ws.Range("E26").Select()
$excel.ActiveWindow.Zoom = 70
